$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the header row (row 1) text ---
$ws.Range("B1").Value = "rural communities"
$ws.Range("C1").Value = "small growers"
$ws.Range("D1").Value = "investor growers"
$ws.Range("E1").Value = "small growers (white area)"
$ws.Range("F1").Value = "investor growers (white area)"
$ws.Range("G1").Value = "municipalities"
# H1 text unchanged ("other dischargers"), only formatting changes

# --- Row height change for the header row ---
$ws.Rows.Item(1).RowHeight = 58.5

# --- Formatting: B1 ---
$ws.Range("B1").WrapText = $true
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Borders.Item(7).LineStyle = 1
$ws.Range("B1").Borders.Item(7).Weight = -4138
$ws.Range("B1").Borders.Item(7).Color = 13421772
$ws.Range("B1").Borders.Item(10).LineStyle = 1
$ws.Range("B1").Borders.Item(10).Weight = -4138
$ws.Range("B1").Borders.Item(10).Color = 13421772
$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("B1").Borders.Item(8).Weight = -4138
$ws.Range("B1").Borders.Item(8).Color = 13421772
$ws.Range("B1").Borders.Item(9).LineStyle = 1
$ws.Range("B1").Borders.Item(9).Weight = -4138
$ws.Range("B1").Borders.Item(9).Color = 0

# --- Formatting: C1 ---
$ws.Range("C1").WrapText = $true
$ws.Range("C1").Font.Name = "Calibri"
$ws.Range("C1").Font.Size = 11
$ws.Range("C1").Borders.Item(7).LineStyle = 1
$ws.Range("C1").Borders.Item(7).Weight = -4138
$ws.Range("C1").Borders.Item(7).Color = 13421772
$ws.Range("C1").Borders.Item(10).LineStyle = 1
$ws.Range("C1").Borders.Item(10).Weight = 2
$ws.Range("C1").Borders.Item(10).Color = 0
$ws.Range("C1").Borders.Item(8).LineStyle = 0
$ws.Range("C1").Borders.Item(9).LineStyle = 1
$ws.Range("C1").Borders.Item(9).Weight = -4138
$ws.Range("C1").Borders.Item(9).Color = 13421772

# --- Formatting: D1:H1 (same style) ---
$ws.Range("D1:H1").WrapText = $true
$ws.Range("D1:H1").Font.Name = "Calibri"
$ws.Range("D1:H1").Font.Size = 11
$ws.Range("D1:H1").Borders.Item(7).LineStyle = 1
$ws.Range("D1:H1").Borders.Item(7).Weight = -4138
$ws.Range("D1:H1").Borders.Item(7).Color = 13421772
$ws.Range("D1:H1").Borders.Item(10).LineStyle = 1
$ws.Range("D1:H1").Borders.Item(10).Weight = 2
$ws.Range("D1:H1").Borders.Item(10).Color = 0
$ws.Range("D1:H1").Borders.Item(8).LineStyle = 1
$ws.Range("D1:H1").Borders.Item(8).Weight = -4138
$ws.Range("D1:H1").Borders.Item(8).Color = 13421772
$ws.Range("D1:H1").Borders.Item(9).LineStyle = 1
$ws.Range("D1:H1").Borders.Item(9).Weight = -4138
$ws.Range("D1:H1").Borders.Item(9).Color = 13421772

# --- sheetView selection change ---
$ws.Range("B1:H1").Select()
